# Adapt column header formatting to respective input file names:
#   "<name>_old" -> "<name>_FV2404"
#   "<name>_new" -> "<name>_FV2410"
# Also wrap the header/data range in an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Determine the used range (header row + data rows) and turn it into a Table.
$lastRow = $ws.Cells(1048576, 1).End(-4162).Row   # xlUp
$lastCol = $headers.Length
$headerRange = $ws.Range($ws.Cells(1, 1), $ws.Cells(1, $lastCol))
$fullRange = $ws.Range($ws.Cells(1, 1), $ws.Cells($lastRow, $lastCol))

$lo = $ws.ListObjects.Add(1, $fullRange, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# Freeze the header row (split below row 1) and keep A2 as the top-left
# cell of the scrollable pane, mirroring the authored sheetView/pane setup.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
